$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D221").Value = 44460
$ws.Range("K221").Value = 600
$ws.Range("L221").Value = 600
$ws.Range("M221").Value = 600
$ws.Range("P221").Value = 600
$ws.Range("D222").Value = 44460
$ws.Range("J222").Value = 1500
$ws.Range("D223").Value = 44414
$ws.Range("J223").Value = 1600
$ws.Range("K223").Value = 700
$ws.Range("L223").Value = 700
$ws.Range("M223").Value = 700
$ws.Range("P223").Value = 700
$ws.Range("D224").Value = 44414
$ws.Range("J224").Value = 1800
$ws.Range("D225").Value = 44165
$ws.Range("J225").Value = 1200
$ws.Range("K225").Value = 600
$ws.Range("L225").Value = 600
$ws.Range("M225").Value = 600
$ws.Range("P225").Value = 600
$ws.Range("D226").Value = 44165
$ws.Range("J226").Value = 1000
$ws.Range("K226").Value = 500
$ws.Range("L226").Value = 500
$ws.Range("M226").Value = 500
$ws.Range("P226").Value = 500
$ws.Range("D227").Value = 44427
$ws.Range("J227").Value = 1000
$ws.Range("D228").Value = 44427
$ws.Range("J228").Value = 1100
$ws.Range("D229").Value = 44172
$ws.Range("J229").Value = 950
$ws.Range("K229").Value = 700
$ws.Range("M229").Value = 700
$ws.Range("P229").Value = 700
$ws.Range("D230").Value = 44172
$ws.Range("J230").Value = 850
$ws.Range("K230").Value = 600
$ws.Range("L230").Value = 600
$ws.Range("M230").Value = 600
$ws.Range("P230").Value = 600
$ws.Range("D231").Value = 44389
$ws.Range("J231").Value = 3380
$ws.Range("K231").Value = 600
$ws.Range("L231").Value = 700
$ws.Range("M231").Value = 656
$ws.Range("P231").Value = 656
$ws.Range("D232").Value = 44389
$ws.Range("I232").Value = "Segunda"
$ws.Range("J232").Value = 1500
$ws.Range("K232").Value = 500
$ws.Range("L232").Value = 500
$ws.Range("M232").Value = 500
$ws.Range("P232").Value = 500
$ws.Range("D233").Value = 44249
$ws.Range("J233").Value = 1200
$ws.Range("D234").Value = 44265
$ws.Range("J234").Value = 1300
$ws.Range("D235").Value = 44343
$ws.Range("J235").Value = 2500
$ws.Range("K235").Value = 800
$ws.Range("L235").Value = 800
$ws.Range("M235").Value = 800
$ws.Range("P235").Value = 800
$ws.Range("D236").Value = 44201
$ws.Range("I236").Value = "Primera"
$ws.Range("J236").Value = 950
$ws.Range("D237").Value = 44280
$ws.Range("J237").Value = 1500
$ws.Range("K237").Value = 1000
$ws.Range("L237").Value = 1000
$ws.Range("M237").Value = 1000
$ws.Range("P237").Value = 1000
$ws.Range("D238").Value = 44280
$ws.Range("J238").Value = 1600
$ws.Range("K238").Value = 800
$ws.Range("L238").Value = 800
$ws.Range("M238").Value = 800
$ws.Range("P238").Value = 800
$ws.Range("D239").Value = 44447
$ws.Range("J239").Value = 1300
$ws.Range("K239").Value = 700
$ws.Range("L239").Value = 700
$ws.Range("M239").Value = 700
$ws.Range("P239").Value = 700
$ws.Range("D240").Value = 44447
$ws.Range("I240").Value = "Segunda"
$ws.Range("J240").Value = 1200
$ws.Range("K240").Value = 500
$ws.Range("L240").Value = 500
$ws.Range("M240").Value = 500
$ws.Range("P240").Value = 500
$ws.Range("D241").Value = 44270
$ws.Range("I241").Value = "Primera"
$ws.Range("J241").Value = 3400
$ws.Range("K241").Value = 800
$ws.Range("L241").Value = 850
$ws.Range("M241").Value = 824
$ws.Range("P241").Value = 824
$ws.Range("D242").Value = 44260
$ws.Range("J242").Value = 750
$ws.Range("K242").Value = 900
$ws.Range("M242").Value = 900
$ws.Range("P242").Value = 900
$ws.Range("D243").Value = 44260
$ws.Range("I243").Value = "Segunda"
$ws.Range("J243").Value = 820
$ws.Range("K243").Value = 700
$ws.Range("L243").Value = 700
$ws.Range("M243").Value = 700
$ws.Range("P243").Value = 700
$ws.Range("D244").Value = 44267
$ws.Range("I244").Value = "Primera"
$ws.Range("J244").Value = 1750
$ws.Range("K244").Value = 800
$ws.Range("L244").Value = 900
$ws.Range("M244").Value = 851
$ws.Range("P244").Value = 851
$ws.Range("D245").Value = 44312
$ws.Range("J245").Value = 950
$ws.Range("K245").Value = 900
$ws.Range("L245").Value = 900
$ws.Range("M245").Value = 900
$ws.Range("P245").Value = 900
$ws.Range("D246").Value = 44312
$ws.Range("J246").Value = 850
$ws.Range("K246").Value = 700
$ws.Range("L246").Value = 700
$ws.Range("M246").Value = 700
$ws.Range("P246").Value = 700
$ws.Range("D247").Value = 44187
$ws.Range("J247").Value = 650
$ws.Range("L247").Value = 600
$ws.Range("M247").Value = 600
$ws.Range("P247").Value = 600
$ws.Range("D248").Value = 44187
$ws.Range("J248").Value = 620
$ws.Range("D249").Value = 44390
$ws.Range("J249").Value = 2450
$ws.Range("K249").Value = 600
$ws.Range("M249").Value = 651
$ws.Range("P249").Value = 651
$ws.Range("D250").Value = 44390
$ws.Range("I250").Value = "Segunda"
$ws.Range("J250").Value = 1700
$ws.Range("K250").Value = 500
$ws.Range("L250").Value = 500
$ws.Range("M250").Value = 500
$ws.Range("P250").Value = 500
$ws.Range("D251").Value = 44386
$ws.Range("I251").Value = "Primera"
$ws.Range("J251").Value = 1200
$ws.Range("D252").Value = 44308
$ws.Range("J252").Value = 2050
$ws.Range("K252").Value = 850
$ws.Range("L252").Value = 900
$ws.Range("M252").Value = 871
$ws.Range("P252").Value = 871
$ws.Range("D253").Value = 44308
$ws.Range("J253").Value = 900
$ws.Range("D254").Value = 44264
$ws.Range("J254").Value = 880
$ws.Range("K254").Value = 800
$ws.Range("L254").Value = 800
$ws.Range("M254").Value = 800
$ws.Range("P254").Value = 800
$ws.Range("D255").Value = 44264
$ws.Range("D256").Value = 44196
$ws.Range("D257").Value = 44196
$ws.Range("J257").Value = 850
$ws.Range("D258").Value = 44301
$ws.Range("K258").Value = 900
$ws.Range("L258").Value = 900
$ws.Range("M258").Value = 900
$ws.Range("P258").Value = 900
$ws.Range("D259").Value = 44301
$ws.Range("I259").Value = "Segunda"
$ws.Range("J259").Value = 950
$ws.Range("K259").Value = 700
$ws.Range("L259").Value = 700
$ws.Range("M259").Value = 700
$ws.Range("P259").Value = 700
$ws.Range("D260").Value = 44251
$ws.Range("I260").Value = "Primera"
$ws.Range("J260").Value = 1200
$ws.Range("K260").Value = 800
$ws.Range("L260").Value = 800
$ws.Range("M260").Value = 800
$ws.Range("P260").Value = 800
$ws.Range("D261").Value = 44243
$ws.Range("J261").Value = 1470
$ws.Range("K261").Value = 650
$ws.Range("L261").Value = 900
$ws.Range("M261").Value = 798
$ws.Range("P261").Value = 798
$ws.Range("D262").Value = 44243
$ws.Range("I262").Value = "Segunda"
$ws.Range("J262").Value = 650
$ws.Range("K262").Value = 700
$ws.Range("L262").Value = 700
$ws.Range("M262").Value = 700
$ws.Range("P262").Value = 700
$ws.Range("D263").Value = 44252
$ws.Range("I263").Value = "Primera"
$ws.Range("J263").Value = 1600
$ws.Range("K263").Value = 800
$ws.Range("L263").Value = 800
$ws.Range("M263").Value = 800
$ws.Range("P263").Value = 800
$ws.Range("D264").Value = 44166
$ws.Range("J264").Value = 980
$ws.Range("L264").Value = 600
$ws.Range("M264").Value = 600
$ws.Range("P264").Value = 600
$ws.Range("D265").Value = 44166
$ws.Range("J265").Value = 970
$ws.Range("K265").Value = 500
$ws.Range("L265").Value = 500
$ws.Range("M265").Value = 500
$ws.Range("P265").Value = 500
$ws.Range("D266").Value = 44168
$ws.Range("J266").Value = 1750
$ws.Range("K266").Value = 600
$ws.Range("L266").Value = 650
$ws.Range("M266").Value = 626
$ws.Range("P266").Value = 626
$ws.Range("D267").Value = 44168
$ws.Range("J267").Value = 900
$ws.Range("K267").Value = 550
$ws.Range("L267").Value = 550
$ws.Range("M267").Value = 550
$ws.Range("P267").Value = 550
$ws.Range("D268").Value = 44369
$ws.Range("J268").Value = 3100
$ws.Range("L268").Value = 750
$ws.Range("M268").Value = 726
$ws.Range("P268").Value = 726
$ws.Range("D269").Value = 44369
$ws.Range("J269").Value = 1500
$ws.Range("K269").Value = 600
$ws.Range("L269").Value = 600
$ws.Range("M269").Value = 600
$ws.Range("P269").Value = 600
$ws.Range("D270").Value = 44433
$ws.Range("J270").Value = 1200
$ws.Range("K270").Value = 700
$ws.Range("L270").Value = 700
$ws.Range("M270").Value = 700
$ws.Range("P270").Value = 700
$ws.Range("D271").Value = 44433
$ws.Range("J271").Value = 900
$ws.Range("K271").Value = 500
$ws.Range("L271").Value = 500
$ws.Range("M271").Value = 500
$ws.Range("P271").Value = 500
$ws.Range("D272").Value = 44221
$ws.Range("J272").Value = 1430
$ws.Range("K272").Value = 850
$ws.Range("L272").Value = 900
$ws.Range("M272").Value = 874
$ws.Range("P272").Value = 874
$ws.Range("D273").Value = 44221
$ws.Range("I273").Value = "Segunda"
$ws.Range("J273").Value = 670
$ws.Range("K273").Value = 700
$ws.Range("L273").Value = 700
$ws.Range("M273").Value = 700
$ws.Range("P273").Value = 700
$ws.Range("D274").Value = 44371
$ws.Range("J274").Value = 2150
$ws.Range("K274").Value = 600
$ws.Range("L274").Value = 700
$ws.Range("M274").Value = 644
$ws.Range("P274").Value = 644
$ws.Range("D275").Value = 44316
$ws.Range("I275").Value = "Primera"
$ws.Range("J275").Value = 1200
$ws.Range("D276").Value = 44279
$ws.Range("J276").Value = 950
$ws.Range("K276").Value = 1000
$ws.Range("L276").Value = 1000
$ws.Range("M276").Value = 1000
$ws.Range("P276").Value = 1000
$ws.Range("D277").Value = 44279
$ws.Range("I277").Value = "Segunda"
$ws.Range("J277").Value = 850
$ws.Range("K277").Value = 800
$ws.Range("L277").Value = 800
$ws.Range("M277").Value = 800
$ws.Range("P277").Value = 800
$ws.Range("D278").Value = 44397
$ws.Range("J278").Value = 1600
$ws.Range("K278").Value = 700
$ws.Range("L278").Value = 700
$ws.Range("M278").Value = 700
$ws.Range("P278").Value = 700
$ws.Range("D279").Value = 44363
$ws.Range("I279").Value = "Primera"
$ws.Range("J279").Value = 2600
$ws.Range("K279").Value = 700
$ws.Range("L279").Value = 700
$ws.Range("M279").Value = 700
$ws.Range("P279").Value = 700
$ws.Range("D280").Value = 44277
$ws.Range("J280").Value = 900
$ws.Range("D281").Value = 44277
$ws.Range("J281").Value = 880
$ws.Range("D282").Value = 44291
$ws.Range("J282").Value = 780
$ws.Range("K282").Value = 1000
$ws.Range("L282").Value = 1000
$ws.Range("M282").Value = 1000
$ws.Range("P282").Value = 1000
$ws.Range("D283").Value = 44291
$ws.Range("J283").Value = 760
$ws.Range("K283").Value = 800
$ws.Range("L283").Value = 800
$ws.Range("M283").Value = 800
$ws.Range("P283").Value = 800
$ws.Range("D284").Value = 44273
$ws.Range("J284").Value = 2400
$ws.Range("K284").Value = 850
$ws.Range("L284").Value = 900
$ws.Range("M284").Value = 875
$ws.Range("P284").Value = 875
$ws.Range("D285").Value = 44273
$ws.Range("J285").Value = 900
$ws.Range("K285").Value = 700
$ws.Range("L285").Value = 700
$ws.Range("M285").Value = 700
$ws.Range("P285").Value = 700
$ws.Range("D286").Value = 44438
$ws.Range("J286").Value = 3800
$ws.Range("K286").Value = 600
$ws.Range("L286").Value = 650
$ws.Range("M286").Value = 625
$ws.Range("P286").Value = 625
$ws.Range("D287").Value = 44438
$ws.Range("I287").Value = "Segunda"
$ws.Range("J287").Value = 1800
$ws.Range("K287").Value = 500
$ws.Range("L287").Value = 500
$ws.Range("M287").Value = 500
$ws.Range("P287").Value = 500
$ws.Range("D288").Value = 44372
$ws.Range("I288").Value = "Primera"
$ws.Range("J288").Value = 2250
$ws.Range("K288").Value = 650
$ws.Range("M288").Value = 679
$ws.Range("P288").Value = 679
$ws.Range("D289").Value = 44286
$ws.Range("K289").Value = 900
$ws.Range("L289").Value = 900
$ws.Range("M289").Value = 900
$ws.Range("P289").Value = 900
$ws.Range("D290").Value = 44286
$ws.Range("J290").Value = 850
$ws.Range("D291").Value = 44209
$ws.Range("J291").Value = 900
$ws.Range("K291").Value = 800
$ws.Range("M291").Value = 800
$ws.Range("P291").Value = 800
$ws.Range("D292").Value = 44209
$ws.Range("J292").Value = 950
$ws.Range("K292").Value = 700
$ws.Range("L292").Value = 700
$ws.Range("M292").Value = 700
$ws.Range("P292").Value = 700
$ws.Range("D293").Value = 44356
$ws.Range("J293").Value = 1790
$ws.Range("K293").Value = 750
$ws.Range("L293").Value = 800
$ws.Range("M293").Value = 774
$ws.Range("P293").Value = 774
$ws.Range("D294").Value = 44356
$ws.Range("J294").Value = 850
$ws.Range("K294").Value = 650
$ws.Range("L294").Value = 650
$ws.Range("M294").Value = 650
$ws.Range("P294").Value = 650
$ws.Range("D295").Value = 44160
$ws.Range("J295").Value = 1750
$ws.Range("K295").Value = 600
$ws.Range("L295").Value = 650
$ws.Range("M295").Value = 624
$ws.Range("P295").Value = 624
$ws.Range("D296").Value = 44160
$ws.Range("J296").Value = 900
$ws.Range("K296").Value = 500
$ws.Range("L296").Value = 500
$ws.Range("M296").Value = 500
$ws.Range("P296").Value = 500
$ws.Range("D297").Value = 44351
$ws.Range("J297").Value = 1150
$ws.Range("K297").Value = 800
$ws.Range("L297").Value = 800
$ws.Range("M297").Value = 800
$ws.Range("P297").Value = 800
$ws.Range("D298").Value = 44351
$ws.Range("J298").Value = 950
$ws.Range("K298").Value = 700
$ws.Range("L298").Value = 700
$ws.Range("M298").Value = 700
$ws.Range("P298").Value = 700
$ws.Range("D299").Value = 44365
$ws.Range("J299").Value = 3300
$ws.Range("K299").Value = 650
$ws.Range("L299").Value = 700
$ws.Range("M299").Value = 677
$ws.Range("P299").Value = 677
$ws.Range("D300").Value = 44365
$ws.Range("J300").Value = 850
$ws.Range("K300").Value = 600
$ws.Range("L300").Value = 600
$ws.Range("M300").Value = 600
$ws.Range("P300").Value = 600
$ws.Range("D301").Value = 44306
$ws.Range("J301").Value = 850
$ws.Range("K301").Value = 900
$ws.Range("L301").Value = 900
$ws.Range("M301").Value = 900
$ws.Range("P301").Value = 900
$ws.Range("D302").Value = 44306
$ws.Range("J302").Value = 950
$ws.Range("K302").Value = 700
$ws.Range("L302").Value = 700
$ws.Range("M302").Value = 700
$ws.Range("P302").Value = 700
$ws.Range("D303").Value = 44215
$ws.Range("J303").Value = 1250
$ws.Range("K303").Value = 750
$ws.Range("L303").Value = 750
$ws.Range("M303").Value = 750
$ws.Range("P303").Value = 750
$ws.Range("D304").Value = 44215
$ws.Range("I304").Value = "Segunda"
$ws.Range("J304").Value = 880
$ws.Range("K304").Value = 650
$ws.Range("L304").Value = 650
$ws.Range("M304").Value = 650
$ws.Range("O304").Value = "Provincia de Quillota"
$ws.Range("P304").Value = 650
$ws.Range("D305").Value = 44175
$ws.Range("I305").Value = "Primera"
$ws.Range("J305").Value = 1750
$ws.Range("K305").Value = 650
$ws.Range("M305").Value = 676
$ws.Range("O305").Value = "Provincia de Quillota"
$ws.Range("P305").Value = 676
$ws.Range("D306").Value = 44357
$ws.Range("J306").Value = 1800
$ws.Range("O306").Value = "Provincia de Santiago"
$ws.Range("D307").Value = 44357
$ws.Range("J307").Value = 1600
$ws.Range("K307").Value = 700
$ws.Range("L307").Value = 700
$ws.Range("M307").Value = 700
$ws.Range("O307").Value = "Provincia de Santiago"
$ws.Range("P307").Value = 700
$ws.Range("D308").Value = 44203
$ws.Range("J308").Value = 850
$ws.Range("K308").Value = 800
$ws.Range("L308").Value = 800
$ws.Range("M308").Value = 800
$ws.Range("P308").Value = 800
$ws.Range("D309").Value = 44203
$ws.Range("I309").Value = "Segunda"
$ws.Range("J309").Value = 900
$ws.Range("K309").Value = 650
$ws.Range("L309").Value = 650
$ws.Range("M309").Value = 650
$ws.Range("P309").Value = 650
$ws.Range("D310").Value = 44162
$ws.Range("J310").Value = 1800
$ws.Range("K310").Value = 700
$ws.Range("L310").Value = 750
$ws.Range("M310").Value = 724
$ws.Range("P310").Value = 724
$ws.Range("D311").Value = 44410
$ws.Range("I311").Value = "Primera"
$ws.Range("J311").Value = 3400
$ws.Range("K311").Value = 600
$ws.Range("L311").Value = 700
$ws.Range("M311").Value = 647
$ws.Range("P311").Value = 647
$ws.Range("D312").Value = 44411
$ws.Range("J312").Value = 3000
$ws.Range("K312").Value = 650
$ws.Range("L312").Value = 700
$ws.Range("M312").Value = 675
$ws.Range("P312").Value = 675
$ws.Range("D313").Value = 44411
$ws.Range("I313").Value = "Segunda"
$ws.Range("J313").Value = 1600
$ws.Range("K313").Value = 500
$ws.Range("L313").Value = 500
$ws.Range("M313").Value = 500
$ws.Range("P313").Value = 500
$ws.Range("D314").Value = 44257
$ws.Range("I314").Value = "Primera"
$ws.Range("K314").Value = 900
$ws.Range("L314").Value = 900
$ws.Range("M314").Value = 900
$ws.Range("P314").Value = 900
$ws.Range("D315").Value = 44244
$ws.Range("J315").Value = 950
$ws.Range("K315").Value = 800
$ws.Range("L315").Value = 800
$ws.Range("M315").Value = 800
$ws.Range("P315").Value = 800
$ws.Range("D316").Value = 44244
$ws.Range("J316").Value = 850
$ws.Range("K316").Value = 700
$ws.Range("L316").Value = 700
$ws.Range("M316").Value = 700
$ws.Range("P316").Value = 700
$ws.Range("D317").Value = 44176
$ws.Range("J317").Value = 2080
$ws.Range("K317").Value = 600
$ws.Range("L317").Value = 650
$ws.Range("M317").Value = 629
$ws.Range("P317").Value = 629
$ws.Range("D318").Value = 44176
$ws.Range("J318").Value = 900
$ws.Range("K318").Value = 550
$ws.Range("L318").Value = 550
$ws.Range("M318").Value = 550
$ws.Range("P318").Value = 550
$ws.Range("D319").Value = 44239
$ws.Range("J319").Value = 900
$ws.Range("K319").Value = 800
$ws.Range("M319").Value = 800
$ws.Range("P319").Value = 800
$ws.Range("D320").Value = 44239
$ws.Range("I320").Value = "Segunda"
$ws.Range("J320").Value = 850
$ws.Range("K320").Value = 700
$ws.Range("L320").Value = 700
$ws.Range("M320").Value = 700
$ws.Range("P320").Value = 700
$ws.Range("D321").Value = 44376
$ws.Range("I321").Value = "Primera"
$ws.Range("J321").Value = 1750
$ws.Range("K321").Value = 700
$ws.Range("M321").Value = 749
$ws.Range("P321").Value = 749
$ws.Range("D322").Value = 44292
$ws.Range("J322").Value = 780
$ws.Range("K322").Value = 1000
$ws.Range("L322").Value = 1000
$ws.Range("M322").Value = 1000
$ws.Range("P322").Value = 1000
$ws.Range("D323").Value = 44292
$ws.Range("J323").Value = 750
$ws.Range("K323").Value = 800
$ws.Range("L323").Value = 800
$ws.Range("M323").Value = 800
$ws.Range("P323").Value = 800
$ws.Range("D324").Value = 44358
$ws.Range("K324").Value = 800
$ws.Range("M324").Value = 800
$ws.Range("P324").Value = 800
$ws.Range("D325").Value = 44358
$ws.Range("J325").Value = 450
$ws.Range("K325").Value = 650
$ws.Range("L325").Value = 650
$ws.Range("M325").Value = 650
$ws.Range("P325").Value = 650
$ws.Range("D326").Value = 44211
$ws.Range("K326").Value = 750
$ws.Range("L326").Value = 800
$ws.Range("M326").Value = 774
$ws.Range("P326").Value = 774
$ws.Range("D327").Value = 44211
$ws.Range("J327").Value = 880
$ws.Range("D328").Value = 44425
$ws.Range("J328").Value = 1800
$ws.Range("K328").Value = 700
$ws.Range("L328").Value = 700
$ws.Range("M328").Value = 700
$ws.Range("P328").Value = 700
$ws.Range("A329").Value = 3
$ws.Range("B329").Value = "Femacal de La Calera"
$ws.Range("C329").Value = "Coquimbo"
$ws.Range("D329").Value = 44425
$ws.Range("E329").Value = 5
$ws.Range("F329").Value = 100112006
$ws.Range("G329").Value = "Repollo"
$ws.Range("H329").Value = "Crespo record"
$ws.Range("I329").Value = "Segunda"
$ws.Range("J329").Value = 900
$ws.Range("K329").Value = 600
$ws.Range("L329").Value = 600
$ws.Range("M329").Value = 600
$ws.Range("N329").Value = "`$/unidad"
$ws.Range("O329").Value = "Provincia de Quillota"
$ws.Range("P329").Value = 600
$ws.Range("Q329").Value = 1
$ws.Range("R329").Value = "Hortaliza"
$ws.Range("D329").NumberFormat = $ws.Range("D220").NumberFormat
$ws.Range("A330").Value = 3
$ws.Range("B330").Value = "Femacal de La Calera"
$ws.Range("C330").Value = "Coquimbo"
$ws.Range("D330").Value = 44323
$ws.Range("E330").Value = 5
$ws.Range("F330").Value = 100112006
$ws.Range("G330").Value = "Repollo"
$ws.Range("H330").Value = "Crespo record"
$ws.Range("I330").Value = "Primera"
$ws.Range("J330").Value = 2500
$ws.Range("K330").Value = 800
$ws.Range("L330").Value = 800
$ws.Range("M330").Value = 800
$ws.Range("N330").Value = "`$/unidad"
$ws.Range("O330").Value = "Provincia de Quillota"
$ws.Range("P330").Value = 800
$ws.Range("Q330").Value = 1
$ws.Range("R330").Value = "Hortaliza"
$ws.Range("D330").NumberFormat = $ws.Range("D220").NumberFormat